$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "num_layers"
$ws.Range("L1").Value = "d_model"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 16

$ws.Columns.Item(11).ColumnWidth = 10

[void]$ws.Range("K1").Select()
